$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the changed cells so numeric-looking / percent-looking
# strings are preserved literally as text (matching original inline-string cells).
$ws.Range("D2:E26").NumberFormat = "@"
$ws.Range("D39:E51").NumberFormat = "@"

$ws.Range("D2").Value = '290.31'
$ws.Range("E2").Value = '-2.82%'
$ws.Range("D3").Value = '30.72'
$ws.Range("E3").Value = '-5.45%'
$ws.Range("D4").Value = '4.918'
$ws.Range("E4").Value = '-2.53%'
$ws.Range("D5").Value = '0.07244'
$ws.Range("E5").Value = '-5.23%'
$ws.Range("D6").Value = '1.808'
$ws.Range("E6").Value = '-10.35%'
$ws.Range("D7").Value = '7.652'
$ws.Range("E7").Value = '-2.81%'
$ws.Range("D8").Value = '3.700'
$ws.Range("E8").Value = '-2.90%'
$ws.Range("D9").Value = '0.9009'
$ws.Range("E9").Value = '-2.65%'
$ws.Range("D10").Value = '0.1685'
$ws.Range("E10").Value = '-3.82%'
$ws.Range("D11").Value = '0.08034'
$ws.Range("E11").Value = '0.21%'
$ws.Range("D12").Value = '0.08076'
$ws.Range("E12").Value = '-5.85%'
$ws.Range("D13").Value = '0.03058'
$ws.Range("E13").Value = '-1.36%'
$ws.Range("D14").Value = '0.1003'
$ws.Range("E14").Value = '0.49%'
$ws.Range("D15").Value = '0.001498'
$ws.Range("E15").Value = '-1.40%'
$ws.Range("D16").Value = '0.005693'
$ws.Range("E16").Value = '-3.62%'
$ws.Range("D17").Value = '3.475'
$ws.Range("D18").Value = '2.077'
$ws.Range("E18").Value = '-3.57%'
$ws.Range("D19").Value = '0.3315'
$ws.Range("E19").Value = '-0.42%'
$ws.Range("D20").Value = '0.1302'
$ws.Range("E20").Value = '-1.75%'
$ws.Range("D21").Value = '3.965'
$ws.Range("E21").Value = '-9.88%'
$ws.Range("E22").Value = '9.53%'
$ws.Range("D23").Value = '0.04524'
$ws.Range("E23").Value = '-0.62%'
$ws.Range("D24").Value = '0.001214'
$ws.Range("E24").Value = '-1.72%'
$ws.Range("D25").Value = '0.004423'
$ws.Range("E25").Value = '7.23%'
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").Value = '3.39%'
$ws.Range("D39").Value = '0.01590'
$ws.Range("E39").Value = '-7.38%'
$ws.Range("E40").Value = '-6.63%'
$ws.Range("D41").Value = '0.007279'
$ws.Range("E41").Value = '-3.10%'
$ws.Range("D43").Value = '0.1314'
$ws.Range("E43").Value = '-2.94%'
$ws.Range("D44").Value = '0.002004'
$ws.Range("E44").Value = '-9.87%'
$ws.Range("D45").Value = '0.009465'
$ws.Range("E45").Value = '-16.65%'
$ws.Range("D46").Value = '0.00005808'
$ws.Range("E46").Value = '-6.06%'
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").Value = '-0.41%'
$ws.Range("D48").Value = '2.254'
$ws.Range("E48").Value = '19.97%'
$ws.Range("E49").Value = '-3.78%'
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '-0.41%'
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '-0.41%'
